$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{"E" = 3; "G" = 179.1580256666667; "H" = 537.4740770000001; "I" = 0.3468013736386751; "J" = 0.3468013736386751; "K" = 3; "M" = 21.09934133333334; "N" = 63.29802400000001; "O" = 0.2917236204149438; "P" = 0.2917236204149438; "Q" = 3780.116336147095; "R" = 34021.04702532385; "S" = 0.1011701522827499; "T" = 0.1011701522827499}
    3 = @{"E" = 3; "G" = 179.1580256666667; "H" = 537.4740770000001; "I" = 0.3468013736386751; "J" = 0.3468013736386751; "K" = 3; "M" = 35.81943766666667; "N" = 107.458313; "O" = 0.4952465516465762; "P" = 0.4952465516465762; "Q" = 6417.339732850234; "R" = 57756.05759565211; "S" = 0.1717521844008497; "T" = 0.1717521844008497}
    4 = @{"E" = 3; "G" = 179.1580256666667; "H" = 537.4740770000001; "I" = 0.3468013736386751; "J" = 0.3468013736386751; "K" = 3; "M" = 15.40769666666667; "N" = 46.22309; "O" = 0.2130298279384801; "P" = 0.2130298279384801; "Q" = 2760.412514870882; "R" = 24843.71263383793; "S" = 0.0738790369550755; "T" = 0.0738790369550755}
    5 = @{"E" = 3; "G" = 239.807332; "H" = 719.421996; "I" = 0.4642019905988459; "J" = 0.4642019905988459; "K" = 3; "M" = 21.09934133333334; "N" = 63.29802400000001; "O" = 0.2917236204149438; "P" = 0.2917236204149438; "Q" = 5059.77675210399; "R" = 45537.99076893591; "S" = 0.135418685301319; "T" = 0.135418685301319}
    6 = @{"E" = 3; "G" = 239.807332; "H" = 719.421996; "I" = 0.4642019905988459; "J" = 0.4642019905988459; "K" = 3; "M" = 35.81943766666667; "N" = 107.458313; "O" = 0.4952465516465762; "P" = 0.4952465516465762; "Q" = 8589.763780583638; "R" = 77307.87402525275; "S" = 0.2298944351115548; "T" = 0.2298944351115548}
    7 = @{"E" = 3; "G" = 239.807332; "H" = 719.421996; "I" = 0.4642019905988459; "J" = 0.4642019905988459; "K" = 3; "M" = 15.40769666666667; "N" = 46.22309; "O" = 0.2130298279384801; "P" = 0.2130298279384801; "Q" = 3694.878629898627; "R" = 33253.90766908764; "S" = 0.09888887018597209; "T" = 0.0988888701859721}
    8 = @{"E" = 3; "G" = 97.63589966666666; "H" = 292.907699; "I" = 0.1889966357624789; "J" = 0.1889966357624789; "K" = 3; "M" = 21.09934133333334; "N" = 63.29802400000001; "O" = 0.2917236204149438; "P" = 0.2917236204149438; "Q" = 2060.053173454086; "R" = 18540.47856108678; "S" = 0.05513478283087479; "T" = 0.0551347828308748}
    9 = @{"E" = 3; "G" = 97.63589966666666; "H" = 292.907699; "I" = 0.1889966357624789; "J" = 0.1889966357624789; "K" = 3; "M" = 35.81943766666667; "N" = 107.458313; "O" = 0.4952465516465762; "P" = 0.4952465516465762; "Q" = 3497.263022139087; "R" = 31475.36719925178; "S" = 0.09359993213417167; "T" = 0.09359993213417168}
    10 = @{"E" = 3; "G" = 97.63589966666666; "H" = 292.907699; "I" = 0.1889966357624789; "J" = 0.1889966357624789; "K" = 3; "M" = 15.40769666666667; "N" = 46.22309; "O" = 0.2130298279384801; "P" = 0.2130298279384801; "Q" = 1504.344325841101; "R" = 13539.09893256991; "S" = 0.04026192079743248; "T" = 0.04026192079743248}
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
